$wb = $excel.ActiveWorkbook

# Work on the "January" sheet and add the list of names
$ws = $wb.Worksheets.Item("January")

$ws.Range("B2").Value = "Nitesh"
$ws.Range("B3").Value = "Gautami"
$ws.Range("B4").Value = "Pratiksha"
$ws.Range("B5").Value = "Pruthvi"

# Select B6 on the January sheet and make it the active sheet/tab
$ws.Activate()
$ws.Range("B6").Select()
